$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(290, 1).Value = "Login Failure with invalid username or password"
$ws.Cells.Item(290, 2).Value = "FAILED"
$ws.Cells.Item(290, 3).Value = "safari"
$ws.Cells.Item(290, 4).NumberFormat = "@"
$ws.Cells.Item(290, 4).Value = "09.11.22"
$ws.Cells.Item(290, 4).ClearFormats()
$ws.Cells.Item(291, 1).Value = "Login Failure with invalid username or password"
$ws.Cells.Item(291, 2).Value = "FAILED"
$ws.Cells.Item(291, 3).Value = "safari"
$ws.Cells.Item(291, 4).NumberFormat = "@"
$ws.Cells.Item(291, 4).Value = "09.11.22"
$ws.Cells.Item(291, 4).ClearFormats()
$ws.Cells.Item(292, 1).Value = "Login Failure with invalid username or password"
$ws.Cells.Item(292, 2).Value = "FAILED"
$ws.Cells.Item(292, 3).Value = "safari"
$ws.Cells.Item(292, 4).NumberFormat = "@"
$ws.Cells.Item(292, 4).Value = "09.11.22"
$ws.Cells.Item(292, 4).ClearFormats()
$ws.Cells.Item(293, 1).Value = "Login Failure with invalid username or password"
$ws.Cells.Item(293, 2).Value = "FAILED"
$ws.Cells.Item(293, 3).Value = "firefox"
$ws.Cells.Item(293, 4).NumberFormat = "@"
$ws.Cells.Item(293, 4).Value = "09.11.22"
$ws.Cells.Item(293, 4).ClearFormats()
$ws.Cells.Item(294, 1).Value = "Login Failure with invalid username or password"
$ws.Cells.Item(294, 2).Value = "FAILED"
$ws.Cells.Item(294, 3).Value = "firefox"
$ws.Cells.Item(294, 4).NumberFormat = "@"
$ws.Cells.Item(294, 4).Value = "09.11.22"
$ws.Cells.Item(294, 4).ClearFormats()
$ws.Cells.Item(295, 1).Value = "Login Failure with invalid username or password"
$ws.Cells.Item(295, 2).Value = "FAILED"
$ws.Cells.Item(295, 3).Value = "firefox"
$ws.Cells.Item(295, 4).NumberFormat = "@"
$ws.Cells.Item(295, 4).Value = "09.11.22"
$ws.Cells.Item(295, 4).ClearFormats()
$ws.Cells.Item(296, 1).Value = "Login Failure with invalid username or password"
$ws.Cells.Item(296, 2).Value = "FAILED"
$ws.Cells.Item(296, 3).Value = "safari"
$ws.Cells.Item(296, 4).NumberFormat = "@"
$ws.Cells.Item(296, 4).Value = "09.11.22"
$ws.Cells.Item(296, 4).ClearFormats()
$ws.Cells.Item(297, 1).Value = "Login Failure with invalid username or password"
$ws.Cells.Item(297, 2).Value = "PASSED"
$ws.Cells.Item(297, 3).Value = "chrome"
$ws.Cells.Item(297, 4).NumberFormat = "@"
$ws.Cells.Item(297, 4).Value = "09.11.22"
$ws.Cells.Item(297, 4).ClearFormats()
$ws.Cells.Item(298, 1).Value = "Login Failure with invalid username or password"
$ws.Cells.Item(298, 2).Value = "FAILED"
$ws.Cells.Item(298, 3).Value = "safari"
$ws.Cells.Item(298, 4).NumberFormat = "@"
$ws.Cells.Item(298, 4).Value = "09.11.22"
$ws.Cells.Item(298, 4).ClearFormats()
$ws.Cells.Item(299, 1).Value = "Login Failure with invalid username or password"
$ws.Cells.Item(299, 2).Value = "FAILED"
$ws.Cells.Item(299, 3).Value = "safari"
$ws.Cells.Item(299, 4).NumberFormat = "@"
$ws.Cells.Item(299, 4).Value = "09.11.22"
$ws.Cells.Item(299, 4).ClearFormats()
$ws.Cells.Item(300, 1).Value = "Login Failure with invalid username or password"
$ws.Cells.Item(300, 2).Value = "FAILED"
$ws.Cells.Item(300, 3).Value = "chrome"
$ws.Cells.Item(300, 4).NumberFormat = "@"
$ws.Cells.Item(300, 4).Value = "09.11.22"
$ws.Cells.Item(300, 4).ClearFormats()
$ws.Cells.Item(301, 1).Value = "Login Failure with invalid username or password"
$ws.Cells.Item(301, 2).Value = "PASSED"
$ws.Cells.Item(301, 3).Value = "chrome"
$ws.Cells.Item(301, 4).NumberFormat = "@"
$ws.Cells.Item(301, 4).Value = "09.11.22"
$ws.Cells.Item(301, 4).ClearFormats()
$ws.Cells.Item(302, 1).Value = "Login Failure with invalid username or password"
$ws.Cells.Item(302, 2).Value = "FAILED"
$ws.Cells.Item(302, 3).Value = "safari"
$ws.Cells.Item(302, 4).NumberFormat = "@"
$ws.Cells.Item(302, 4).Value = "09.11.22"
$ws.Cells.Item(302, 4).ClearFormats()
$ws.Cells.Item(303, 1).Value = "Login Failure with invalid username or password"
$ws.Cells.Item(303, 2).Value = "FAILED"
$ws.Cells.Item(303, 3).Value = "safari"
$ws.Cells.Item(303, 4).NumberFormat = "@"
$ws.Cells.Item(303, 4).Value = "09.11.22"
$ws.Cells.Item(303, 4).ClearFormats()
$ws.Cells.Item(304, 1).Value = "Login Failure with invalid username or password"
$ws.Cells.Item(304, 2).Value = "FAILED"
$ws.Cells.Item(304, 3).Value = "safari"
$ws.Cells.Item(304, 4).NumberFormat = "@"
$ws.Cells.Item(304, 4).Value = "09.11.22"
$ws.Cells.Item(304, 4).ClearFormats()
$ws.Cells.Item(305, 1).Value = "Login Failure with invalid username or password"
$ws.Cells.Item(305, 2).Value = "PASSED"
$ws.Cells.Item(305, 3).Value = "chrome"
$ws.Cells.Item(305, 4).NumberFormat = "@"
$ws.Cells.Item(305, 4).Value = "09.11.22"
$ws.Cells.Item(305, 4).ClearFormats()
$ws.Cells.Item(306, 1).Value = "Login Failure with invalid username or password"
$ws.Cells.Item(306, 2).Value = "FAILED"
$ws.Cells.Item(306, 3).Value = "safari"
$ws.Cells.Item(306, 4).NumberFormat = "@"
$ws.Cells.Item(306, 4).Value = "09.11.22"
$ws.Cells.Item(306, 4).ClearFormats()
$ws.Cells.Item(307, 1).Value = "Login Failure with invalid username or password"
$ws.Cells.Item(307, 2).Value = "PASSED"
$ws.Cells.Item(307, 3).Value = "chrome"
$ws.Cells.Item(307, 4).NumberFormat = "@"
$ws.Cells.Item(307, 4).Value = "09.11.22"
$ws.Cells.Item(307, 4).ClearFormats()
$ws.Cells.Item(308, 1).Value = "Login Failure with invalid username or password"
$ws.Cells.Item(308, 2).Value = "PASSED"
$ws.Cells.Item(308, 3).Value = "chrome"
$ws.Cells.Item(308, 4).NumberFormat = "@"
$ws.Cells.Item(308, 4).Value = "09.11.22"
$ws.Cells.Item(308, 4).ClearFormats()
$ws.Cells.Item(309, 1).Value = "Login Failure with invalid username or password"
$ws.Cells.Item(309, 2).Value = "FAILED"
$ws.Cells.Item(309, 3).Value = "safari"
$ws.Cells.Item(309, 4).NumberFormat = "@"
$ws.Cells.Item(309, 4).Value = "09.11.22"
$ws.Cells.Item(309, 4).ClearFormats()
$ws.Cells.Item(310, 1).Value = "Login Failure with invalid username or password"
$ws.Cells.Item(310, 2).Value = "PASSED"
$ws.Cells.Item(310, 3).Value = "chrome"
$ws.Cells.Item(310, 4).NumberFormat = "@"
$ws.Cells.Item(310, 4).Value = "09.11.22"
$ws.Cells.Item(310, 4).ClearFormats()
$ws.Cells.Item(311, 1).Value = "Login Failure with invalid username or password"
$ws.Cells.Item(311, 2).Value = "PASSED"
$ws.Cells.Item(311, 3).Value = "firefox"
$ws.Cells.Item(311, 4).NumberFormat = "@"
$ws.Cells.Item(311, 4).Value = "09.11.22"
$ws.Cells.Item(311, 4).ClearFormats()
$ws.Cells.Item(312, 1).Value = "Login Failure with invalid username or password"
$ws.Cells.Item(312, 2).Value = "PASSED"
$ws.Cells.Item(312, 3).Value = "chrome"
$ws.Cells.Item(312, 4).NumberFormat = "@"
$ws.Cells.Item(312, 4).Value = "09.11.22"
$ws.Cells.Item(312, 4).ClearFormats()
$ws.Cells.Item(313, 1).Value = "Login Failure with invalid username or password"
$ws.Cells.Item(313, 2).Value = "FAILED"
$ws.Cells.Item(313, 3).Value = "firefox"
$ws.Cells.Item(313, 4).NumberFormat = "@"
$ws.Cells.Item(313, 4).Value = "09.11.22"
$ws.Cells.Item(313, 4).ClearFormats()
$ws.Cells.Item(314, 1).Value = "Login Failure with invalid username or password"
$ws.Cells.Item(314, 2).Value = "PASSED"
$ws.Cells.Item(314, 3).Value = "chrome"
$ws.Cells.Item(314, 4).NumberFormat = "@"
$ws.Cells.Item(314, 4).Value = "09.11.22"
$ws.Cells.Item(314, 4).ClearFormats()
$ws.Cells.Item(315, 1).Value = "Login with valid username and password"
$ws.Cells.Item(315, 2).Value = "FAILED"
$ws.Cells.Item(315, 3).Value = "chrome"
$ws.Cells.Item(315, 4).NumberFormat = "@"
$ws.Cells.Item(315, 4).Value = "29.12.22"
$ws.Cells.Item(315, 4).ClearFormats()
$ws.Cells.Item(316, 1).Value = "Login with one invalid pair of (username and password)"
$ws.Cells.Item(316, 2).Value = "FAILED"
$ws.Cells.Item(316, 3).Value = "chrome"
$ws.Cells.Item(316, 4).NumberFormat = "@"
$ws.Cells.Item(316, 4).Value = "29.12.22"
$ws.Cells.Item(316, 4).ClearFormats()
$ws.Cells.Item(317, 1).Value = "Login with one invalid pair of (username and password)"
$ws.Cells.Item(317, 2).Value = "FAILED"
$ws.Cells.Item(317, 3).Value = "chrome"
$ws.Cells.Item(317, 4).NumberFormat = "@"
$ws.Cells.Item(317, 4).Value = "29.12.22"
$ws.Cells.Item(317, 4).ClearFormats()
$ws.Cells.Item(318, 1).Value = "Login with one invalid pair of (username and password)"
$ws.Cells.Item(318, 2).Value = "FAILED"
$ws.Cells.Item(318, 3).Value = "chrome"
$ws.Cells.Item(318, 4).NumberFormat = "@"
$ws.Cells.Item(318, 4).Value = "29.12.22"
$ws.Cells.Item(318, 4).ClearFormats()
$ws.Cells.Item(319, 1).Value = "Login with valid username and password"
$ws.Cells.Item(319, 2).Value = "PASSED"
$ws.Cells.Item(319, 3).Value = "chrome"
$ws.Cells.Item(319, 4).NumberFormat = "@"
$ws.Cells.Item(319, 4).Value = "29.12.22"
$ws.Cells.Item(319, 4).ClearFormats()
$ws.Cells.Item(320, 1).Value = "Login with one invalid pair of (username and password)"
$ws.Cells.Item(320, 2).Value = "PASSED"
$ws.Cells.Item(320, 3).Value = "chrome"
$ws.Cells.Item(320, 4).NumberFormat = "@"
$ws.Cells.Item(320, 4).Value = "29.12.22"
$ws.Cells.Item(320, 4).ClearFormats()
$ws.Cells.Item(321, 1).Value = "Login with one invalid pair of (username and password)"
$ws.Cells.Item(321, 2).Value = "PASSED"
$ws.Cells.Item(321, 3).Value = "chrome"
$ws.Cells.Item(321, 4).NumberFormat = "@"
$ws.Cells.Item(321, 4).Value = "29.12.22"
$ws.Cells.Item(321, 4).ClearFormats()
$ws.Cells.Item(322, 1).Value = "Login with one invalid pair of (username and password)"
$ws.Cells.Item(322, 2).Value = "PASSED"
$ws.Cells.Item(322, 3).Value = "chrome"
$ws.Cells.Item(322, 4).NumberFormat = "@"
$ws.Cells.Item(322, 4).Value = "29.12.22"
$ws.Cells.Item(322, 4).ClearFormats()
$ws.Cells.Item(323, 1).Value = "Login Failure with invalid username or password"
$ws.Cells.Item(323, 2).Value = "PASSED"
$ws.Cells.Item(323, 3).Value = "chrome"
$ws.Cells.Item(323, 4).NumberFormat = "@"
$ws.Cells.Item(323, 4).Value = "29.12.22"
$ws.Cells.Item(323, 4).ClearFormats()
$ws.Cells.Item(324, 1).Value = "Login Failure with invalid username or password"
$ws.Cells.Item(324, 2).Value = "PASSED"
$ws.Cells.Item(324, 3).Value = "chrome"
$ws.Cells.Item(324, 4).NumberFormat = "@"
$ws.Cells.Item(324, 4).Value = "29.12.22"
$ws.Cells.Item(324, 4).ClearFormats()
$ws.Cells.Item(325, 1).Value = "Login Failure with invalid username or password"
$ws.Cells.Item(325, 2).Value = "PASSED"
$ws.Cells.Item(325, 3).Value = "chrome"
$ws.Cells.Item(325, 4).NumberFormat = "@"
$ws.Cells.Item(325, 4).Value = "29.12.22"
$ws.Cells.Item(325, 4).ClearFormats()
$ws.Cells.Item(326, 1).Value = "Login Failure with invalid username or password"
$ws.Cells.Item(326, 2).Value = "PASSED"
$ws.Cells.Item(326, 3).Value = "chrome"
$ws.Cells.Item(326, 4).NumberFormat = "@"
$ws.Cells.Item(326, 4).Value = "29.12.22"
$ws.Cells.Item(326, 4).ClearFormats()
$ws.Cells.Item(327, 1).Value = "Login Failure with invalid username or password"
$ws.Cells.Item(327, 2).Value = "FAILED"
$ws.Cells.Item(327, 3).Value = "firefox"
$ws.Cells.Item(327, 4).NumberFormat = "@"
$ws.Cells.Item(327, 4).Value = "29.12.22"
$ws.Cells.Item(327, 4).ClearFormats()
$ws.Cells.Item(328, 1).Value = "Login Failure with invalid username or password"
$ws.Cells.Item(328, 2).Value = "PASSED"
$ws.Cells.Item(328, 3).Value = "chrome"
$ws.Cells.Item(328, 4).NumberFormat = "@"
$ws.Cells.Item(328, 4).Value = "29.12.22"
$ws.Cells.Item(328, 4).ClearFormats()
$ws.Cells.Item(329, 1).Value = "Login Failure with invalid username or password"
$ws.Cells.Item(329, 2).Value = "FAILED"
$ws.Cells.Item(329, 3).Value = "firefox"
$ws.Cells.Item(329, 4).NumberFormat = "@"
$ws.Cells.Item(329, 4).Value = "29.12.22"
$ws.Cells.Item(329, 4).ClearFormats()
$ws.Cells.Item(330, 1).Value = "Login Failure with invalid username or password"
$ws.Cells.Item(330, 2).Value = "PASSED"
$ws.Cells.Item(330, 3).Value = "chrome"
$ws.Cells.Item(330, 4).NumberFormat = "@"
$ws.Cells.Item(330, 4).Value = "29.12.22"
$ws.Cells.Item(330, 4).ClearFormats()
$ws.Cells.Item(331, 1).Value = "Login Failure with invalid username or password"
$ws.Cells.Item(331, 2).Value = "PASSED"
$ws.Cells.Item(331, 3).Value = "chrome"
$ws.Cells.Item(331, 4).NumberFormat = "@"
$ws.Cells.Item(331, 4).Value = "29.12.22"
$ws.Cells.Item(331, 4).ClearFormats()
$ws.Cells.Item(332, 1).Value = "Login Failure with invalid username or password"
$ws.Cells.Item(332, 2).Value = "PASSED"
$ws.Cells.Item(332, 3).Value = "firefox"
$ws.Cells.Item(332, 4).NumberFormat = "@"
$ws.Cells.Item(332, 4).Value = "29.12.22"
$ws.Cells.Item(332, 4).ClearFormats()
$ws.Cells.Item(333, 1).Value = "Login Failure with invalid username or password"
$ws.Cells.Item(333, 2).Value = "PASSED"
$ws.Cells.Item(333, 3).Value = "chrome"
$ws.Cells.Item(333, 4).NumberFormat = "@"
$ws.Cells.Item(333, 4).Value = "29.12.22"
$ws.Cells.Item(333, 4).ClearFormats()
$ws.Cells.Item(334, 1).Value = "Login Failure with invalid username or password"
$ws.Cells.Item(334, 2).Value = "FAILED"
$ws.Cells.Item(334, 3).Value = "firefox"
$ws.Cells.Item(334, 4).NumberFormat = "@"
$ws.Cells.Item(334, 4).Value = "29.12.22"
$ws.Cells.Item(334, 4).ClearFormats()
$ws.Cells.Item(335, 1).Value = "Login Failure with invalid username or password"
$ws.Cells.Item(335, 2).Value = "FAILED"
$ws.Cells.Item(335, 3).Value = "firefox"
$ws.Cells.Item(335, 4).NumberFormat = "@"
$ws.Cells.Item(335, 4).Value = "29.12.22"
$ws.Cells.Item(335, 4).ClearFormats()
$ws.Cells.Item(336, 1).Value = "Login Failure with invalid username or password"
$ws.Cells.Item(336, 2).Value = "PASSED"
$ws.Cells.Item(336, 3).Value = "chrome"
$ws.Cells.Item(336, 4).NumberFormat = "@"
$ws.Cells.Item(336, 4).Value = "29.12.22"
$ws.Cells.Item(336, 4).ClearFormats()
$ws.Cells.Item(337, 1).Value = "Login with valid username and password"
$ws.Cells.Item(337, 2).Value = "PASSED"
$ws.Cells.Item(337, 3).Value = "chrome"
$ws.Cells.Item(337, 4).NumberFormat = "@"
$ws.Cells.Item(337, 4).Value = "29.12.22"
$ws.Cells.Item(337, 4).ClearFormats()
$ws.Cells.Item(338, 1).Value = "Login with valid username and password"
$ws.Cells.Item(338, 2).Value = "FAILED"
$ws.Cells.Item(338, 3).Value = "firefox"
$ws.Cells.Item(338, 4).NumberFormat = "@"
$ws.Cells.Item(338, 4).Value = "29.12.22"
$ws.Cells.Item(338, 4).ClearFormats()
$ws.Cells.Item(339, 1).Value = "Login Failure with invalid username or password"
$ws.Cells.Item(339, 2).Value = "PASSED"
$ws.Cells.Item(339, 3).Value = "chrome"
$ws.Cells.Item(339, 4).NumberFormat = "@"
$ws.Cells.Item(339, 4).Value = "29.12.22"
$ws.Cells.Item(339, 4).ClearFormats()
$ws.Cells.Item(340, 1).Value = "Login Failure with invalid username or password"
$ws.Cells.Item(340, 2).Value = "PASSED"
$ws.Cells.Item(340, 3).Value = "firefox"
$ws.Cells.Item(340, 4).NumberFormat = "@"
$ws.Cells.Item(340, 4).Value = "29.12.22"
$ws.Cells.Item(340, 4).ClearFormats()
$ws.Cells.Item(341, 1).Value = "Login Failure with invalid username or password"
$ws.Cells.Item(341, 2).Value = "FAILED"
$ws.Cells.Item(341, 3).Value = "chrome"
$ws.Cells.Item(341, 4).NumberFormat = "@"
$ws.Cells.Item(341, 4).Value = "29.12.22"
$ws.Cells.Item(341, 4).ClearFormats()
$ws.Cells.Item(342, 1).Value = "Login Failure with invalid username or password"
$ws.Cells.Item(342, 2).Value = "FAILED"
$ws.Cells.Item(342, 3).Value = "firefox"
$ws.Cells.Item(342, 4).NumberFormat = "@"
$ws.Cells.Item(342, 4).Value = "29.12.22"
$ws.Cells.Item(342, 4).ClearFormats()
$ws.Cells.Item(343, 1).Value = "Login Failure with invalid username or password"
$ws.Cells.Item(343, 2).Value = "PASSED"
$ws.Cells.Item(343, 3).Value = "chrome"
$ws.Cells.Item(343, 4).NumberFormat = "@"
$ws.Cells.Item(343, 4).Value = "29.12.22"
$ws.Cells.Item(343, 4).ClearFormats()
$ws.Cells.Item(344, 1).Value = "Login Failure with invalid username or password"
$ws.Cells.Item(344, 2).Value = "FAILED"
$ws.Cells.Item(344, 3).Value = "firefox"
$ws.Cells.Item(344, 4).NumberFormat = "@"
$ws.Cells.Item(344, 4).Value = "29.12.22"
$ws.Cells.Item(344, 4).ClearFormats()
$ws.Cells.Item(345, 1).Value = "Login with valid username and password"
$ws.Cells.Item(345, 2).Value = "PASSED"
$ws.Cells.Item(345, 3).Value = "chrome"
$ws.Cells.Item(345, 4).NumberFormat = "@"
$ws.Cells.Item(345, 4).Value = "29.12.22"
$ws.Cells.Item(345, 4).ClearFormats()
$ws.Cells.Item(346, 1).Value = "Login with valid username and password"
$ws.Cells.Item(346, 2).Value = "FAILED"
$ws.Cells.Item(346, 3).Value = "firefox"
$ws.Cells.Item(346, 4).NumberFormat = "@"
$ws.Cells.Item(346, 4).Value = "29.12.22"
$ws.Cells.Item(346, 4).ClearFormats()
$ws.Cells.Item(347, 1).Value = "Login with valid username and password"
$ws.Cells.Item(347, 2).Value = "PASSED"
$ws.Cells.Item(347, 3).Value = "chrome"
$ws.Cells.Item(347, 4).NumberFormat = "@"
$ws.Cells.Item(347, 4).Value = "29.12.22"
$ws.Cells.Item(347, 4).ClearFormats()
$ws.Cells.Item(348, 1).Value = "Login with one invalid pair of (username and password)"
$ws.Cells.Item(348, 2).Value = "PASSED"
$ws.Cells.Item(348, 3).Value = "chrome"
$ws.Cells.Item(348, 4).NumberFormat = "@"
$ws.Cells.Item(348, 4).Value = "29.12.22"
$ws.Cells.Item(348, 4).ClearFormats()
$ws.Cells.Item(349, 1).Value = "Login with one invalid pair of (username and password)"
$ws.Cells.Item(349, 2).Value = "PASSED"
$ws.Cells.Item(349, 3).Value = "chrome"
$ws.Cells.Item(349, 4).NumberFormat = "@"
$ws.Cells.Item(349, 4).Value = "29.12.22"
$ws.Cells.Item(349, 4).ClearFormats()
$ws.Cells.Item(350, 1).Value = "Login with one invalid pair of (username and password)"
$ws.Cells.Item(350, 2).Value = "PASSED"
$ws.Cells.Item(350, 3).Value = "chrome"
$ws.Cells.Item(350, 4).NumberFormat = "@"
$ws.Cells.Item(350, 4).Value = "29.12.22"
$ws.Cells.Item(350, 4).ClearFormats()
$ws.Cells.Item(351, 1).Value = "Login with valid username and password"
$ws.Cells.Item(351, 2).Value = "PASSED"
$ws.Cells.Item(351, 3).Value = "chrome"
$ws.Cells.Item(351, 4).NumberFormat = "@"
$ws.Cells.Item(351, 4).Value = "29.12.22"
$ws.Cells.Item(351, 4).ClearFormats()
$ws.Cells.Item(352, 1).Value = "Login Failure with invalid username or password"
$ws.Cells.Item(352, 2).Value = "FAILED"
$ws.Cells.Item(352, 3).Value = "chrome"
$ws.Cells.Item(352, 4).NumberFormat = "@"
$ws.Cells.Item(352, 4).Value = "29.12.22"
$ws.Cells.Item(352, 4).ClearFormats()
$ws.Cells.Item(353, 1).Value = "Login Failure with invalid username or password"
$ws.Cells.Item(353, 2).Value = "PASSED"
$ws.Cells.Item(353, 3).Value = "chrome"
$ws.Cells.Item(353, 4).NumberFormat = "@"
$ws.Cells.Item(353, 4).Value = "29.12.22"
$ws.Cells.Item(353, 4).ClearFormats()
$ws.Cells.Item(354, 1).Value = "Login Failure with invalid username or password"
$ws.Cells.Item(354, 2).Value = "PASSED"
$ws.Cells.Item(354, 3).Value = "chrome"
$ws.Cells.Item(354, 4).NumberFormat = "@"
$ws.Cells.Item(354, 4).Value = "29.12.22"
$ws.Cells.Item(354, 4).ClearFormats()
$ws.Cells.Item(355, 1).Value = "Login with valid username and password"
$ws.Cells.Item(355, 2).Value = "FAILED"
$ws.Cells.Item(355, 3).Value = "chrome"
$ws.Cells.Item(355, 4).NumberFormat = "@"
$ws.Cells.Item(355, 4).Value = "29.12.22"
$ws.Cells.Item(355, 4).ClearFormats()
$ws.Cells.Item(356, 1).Value = "Login Failure with invalid username or password"
$ws.Cells.Item(356, 2).Value = "FAILED"
$ws.Cells.Item(356, 3).Value = "chrome"
$ws.Cells.Item(356, 4).NumberFormat = "@"
$ws.Cells.Item(356, 4).Value = "29.12.22"
$ws.Cells.Item(356, 4).ClearFormats()
$ws.Cells.Item(357, 1).Value = "Login Failure with invalid username or password"
$ws.Cells.Item(357, 2).Value = "FAILED"
$ws.Cells.Item(357, 3).Value = "chrome"
$ws.Cells.Item(357, 4).NumberFormat = "@"
$ws.Cells.Item(357, 4).Value = "29.12.22"
$ws.Cells.Item(357, 4).ClearFormats()
$ws.Cells.Item(358, 1).Value = "Login Failure with invalid username or password"
$ws.Cells.Item(358, 2).Value = "FAILED"
$ws.Cells.Item(358, 3).Value = "chrome"
$ws.Cells.Item(358, 4).NumberFormat = "@"
$ws.Cells.Item(358, 4).Value = "29.12.22"
$ws.Cells.Item(358, 4).ClearFormats()
